$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format columns D:G for the data rows as Text so that values such as
# "4-2-2023", "1.73%" and long decimals are stored literally instead of
# being auto-converted by Excel into dates / numbers / percentages.
$ws.Range("D2:G51").NumberFormat = "@"

$data = @{
    2 = @{ D="332.91"; E="1.73%"; F="4-2-2023"; G="1" }
    3 = @{ D="41.21"; E="3.75%"; F="4-2-2023"; G="1" }
    4 = @{ D="5.745"; E="-2.85%"; F="4-2-2023"; G="1" }
    5 = @{ D="0.08213"; E="2.31%"; F="4-2-2023"; G="1" }
    6 = @{ D="2.033"; E="6.06%"; F="4-2-2023"; G="1" }
    7 = @{ D="8.776"; E="1.04%"; F="4-2-2023"; G="1" }
    8 = @{ D="4.536"; E="-1.41%"; F="4-2-2023"; G="1" }
    9 = @{ D="2.926"; E="-0.50%"; F="4-2-2023"; G="1" }
    10 = @{ D="0.9225"; E="-1.06%"; F="4-2-2023"; G="1" }
    11 = @{ D="0.1253"; E="1.31%"; F="4-2-2023"; G="1" }
    12 = @{ D="0.1965"; E="-0.12%"; F="4-2-2023"; G="1" }
    13 = @{ D="8.339"; E="-5.02%"; F="4-2-2023"; G="1" }
    14 = @{ D="0.09403"; E="2.27%"; F="4-2-2023"; G="1" }
    15 = @{ D="0.03671"; E="6.22%"; F="4-2-2023"; G="1" }
    16 = @{ D="0.1054"; E="9.55%"; F="4-2-2023"; G="1" }
    17 = @{ D="0.001306"; E="1.05%"; F="4-2-2023"; G="1" }
    18 = @{ D="0.006259"; E="2.39%"; F="4-2-2023"; G="1" }
    19 = @{ D="3.387"; E="1.53%"; F="4-2-2023"; G="1" }
    20 = @{ D="0.3486"; E="-1.48%"; F="4-2-2023"; G="1" }
    21 = @{ D="0.1418"; E="-0.81%"; F="4-2-2023"; G="1" }
    22 = @{ E="10.00%"; F="4-2-2023"; G="1" }
    23 = @{ D="0.04429"; E="-0.20%"; F="4-2-2023"; G="1" }
    24 = @{ D="0.001270"; E="0.87%"; F="4-2-2023"; G="1" }
    25 = @{ D="0.004309"; E="-1.42%"; F="4-2-2023"; G="1" }
    26 = @{ E="8.80%"; F="4-2-2023"; G="1" }
    27 = @{ F="4-2-2023"; G="1" }
    28 = @{ F="4-2-2023"; G="1" }
    29 = @{ F="4-2-2023"; G="1" }
    30 = @{ F="4-2-2023"; G="1" }
    31 = @{ F="4-2-2023"; G="1" }
    32 = @{ F="4-2-2023"; G="1" }
    33 = @{ F="4-2-2023"; G="1" }
    34 = @{ F="4-2-2023"; G="1" }
    35 = @{ F="4-2-2023"; G="1" }
    36 = @{ F="4-2-2023"; G="1" }
    37 = @{ F="4-2-2023"; G="1" }
    38 = @{ F="4-2-2023"; G="1" }
    39 = @{ D="0.02820"; E="16.27%"; F="4-2-2023"; G="1" }
    40 = @{ D="0.05490"; E="5.06%"; F="4-2-2023"; G="1" }
    41 = @{ D="0.007677"; E="3.11%"; F="4-2-2023"; G="1" }
    42 = @{ D="0.009964"; E="13.73%"; F="4-2-2023"; G="1" }
    43 = @{ D="0.1424"; E="1.08%"; F="4-2-2023"; G="1" }
    44 = @{ D="0.002122"; E="0.03%"; F="4-2-2023"; G="1" }
    45 = @{ D="0.01191"; E="23.49%"; F="4-2-2023"; G="1" }
    46 = @{ D="0.00006822"; E="2.09%"; F="4-2-2023"; G="1" }
    47 = @{ D="0.00000000751"; E="0.02%"; F="4-2-2023"; G="1" }
    48 = @{ D="0.002281"; E="60.35%"; F="4-2-2023"; G="1" }
    49 = @{ D="0.003062"; E="2.07%"; F="4-2-2023"; G="1" }
    50 = @{ D="0.00002102"; E="0.02%"; F="4-2-2023"; G="1" }
    51 = @{ D="0.0002002"; E="0.02%"; F="4-2-2023"; G="1" }
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$r").Value = $rowVals[$col]
    }
}

Write-Host "Applied updates to" $data.Keys.Count "rows"
